$d = $word.ActiveDocument

# --- Locate the paragraph that holds "Time.deltaTime: ... (update) y otra." ---
$para = $d.Paragraphs(3)
$pr = $para.Range

# --- Find the "update" word inside the parenthetical and split it into its own run ---
$fUpdate = $pr.Duplicate
$fUpdate.Find.ClearFormatting()
$fUpdate.Find.Execute("update") | Out-Null
$updateStart = $fUpdate.Start
$updateEnd = $fUpdate.End

# Use a throwaway bookmark at the two boundaries to force Word to split the
# surrounding run into three runs: "...(" / "update" / ")..." without leaving
# any residual character formatting behind.
$rUpdate = $d.Range($updateStart, $updateEnd)
$d.Bookmarks.Add("ZZZsplit1", $rUpdate) | Out-Null
$d.Bookmarks("ZZZsplit1").Delete()

# --- Move the ")" that trails "update" so it precedes " y otra." instead ---
$rParen = $d.Range($updateEnd, $updateEnd + 1)
if ($rParen.Text -eq ")") {
    $rParen.Delete()
    $rTail = $d.Range($updateEnd, $updateEnd)
    $rTail.InsertBefore(")")
}

# --- Split the paragraph right after ") y otra." so the bookmark ends up in its own paragraph ---
$fTail = $pr.Duplicate
$fTail.Find.ClearFormatting()
$fTail.Find.Execute(") y otra.") | Out-Null
$tailEnd = $fTail.End

$rBreak = $d.Range($tailEnd, $tailEnd)
$rBreak.InsertParagraphAfter()

# --- Fill in the new second paragraph with the Time.timeScale entry ---
$newPara = $d.Paragraphs(4)
$newPara.Range.Text = "Time.timeScale: Controla que tan rápido se comporta el juego respecto al tiempo real. 0 es pausado y 1 es tiempo normal."

# --- Re-home the _GoBack bookmark at the very end of the new paragraph ---
$d.Bookmarks("_GoBack").Delete()
$endOfNewPara = $d.Paragraphs(4).Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($endOfNewPara, $endOfNewPara)) | Out-Null
